$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2:K42").Value = "E7420L"
$ws.Range("N6").Select()
